$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new zero-valued cells for F3:G5 (longitudElectrica / Frecuencia columns)
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0

# Update the active selection to match the saved view state
$ws.Range("H17").Select()
